$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Static Object")

# Insert two blank rows to push "coin" and "cloud" blocks down,
# matching the new gaps used to host the "question block" data further below.
$ws.Rows("5:5").Insert() | Out-Null
$ws.Rows("11:11").Insert() | Out-Null

# Mark a section-separator cell (italic style, no value) above the new data.
$ws.Range("A17").Font.Italic = $true

# "question block" / "not hit" state block
$ws.Range("A18").Value = "question block"
$ws.Range("B18").Value = "not hit"
$ws.Range("C18").Value = "not hit"
$ws.Range("D18").Value = 1030001
$ws.Range("E18").Value = "not hit"
$ws.Range("F18").Value = 1030000

$ws.Range("D19").Value = 1030002
$ws.Range("D20").Value = 1030003

# "hit" state block
$ws.Range("B22").Value = "hit"
$ws.Range("C22").Value = "hit "
$ws.Range("D22").Value = 1031001
$ws.Range("E22").Value = "hit"
$ws.Range("F22").Value = 1031000

# Make "Static Object" the active/visible tab, with the new selection.
$ws.Activate()
$ws.Range("G24").Select() | Out-Null
